$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move "applebees / nonfiction / 1" row to the bottom, shift everything else up.
$ws.Range("A2").Value = "dragon ball"
$ws.Range("B2").Value = "manga"
$ws.Range("C2").Value = "1,2,3,4,5,6,7"

$ws.Range("A3").Value = "naruto"
$ws.Range("B3").Value = "manga"
$ws.Range("C3").Value = "2"

$ws.Range("A4").Value = "d"
$ws.Range("B4").Value = "manga"
$ws.Range("C4").Value = "2"

$ws.Range("A5").Value = "dragonbody"
$ws.Range("B5").Value = "manga"
$ws.Range("C5").Value = "2"

$ws.Range("A6").Value = "applebees"
$ws.Range("B6").Value = "nonfiction"
$ws.Range("C6").Value = "1"
